$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "H 72" record (row 2) entirely; all following rows shift up by one.
$ws.Rows.Item(2).Delete()
